# Refresh the cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# figures on Sheet1 to the latest scraped snapshot, as produced by the
# scheduled GitHub Actions job.
#
# Column D cells sometimes hold numeric-looking text (e.g. "1.00", "13.10")
# that must stay as literal text (matching the original report's
# formatting) instead of being auto-coerced into numbers by Excel, which
# would silently drop trailing zeros / introduce floating point noise.
# Setting NumberFormat to "@" (Text) before assigning the value keeps it
# as an exact string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.629.26'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.773.40'
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '355.46'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.97'
$ws.Range("E6").Value = '  -2.70%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.551'
$ws.Range("E7").Value = '  -3.42%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.583'
$ws.Range("E9").Value = '  -2.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.41'
$ws.Range("E10").Value = '  -3.44%  '
$ws.Range("E11").Value = '  +2.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0841'
$ws.Range("E12").Value = '  -2.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.42'
$ws.Range("E13").Value = '  -2.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.55'
$ws.Range("E14").Value = '  -2.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.208.14'
$ws.Range("E15").Value = '  -1.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.764.15'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.927'
$ws.Range("E17").Value = '  -0.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.608.04'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.43'
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.09'
$ws.Range("E20").Value = '  -3.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.10'
$ws.Range("E21").Value = '  -2.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0965'
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.92'
$ws.Range("E23").Value = '  -1.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.48'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.73'
$ws.Range("E25").Value = '  -3.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.29'
$ws.Range("E26").Value = '  -2.58%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.164'
$ws.Range("E28").Value = '  +14.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.15'
$ws.Range("E29").Value = '  -1.70%  '
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.10'
$ws.Range("E31").Value = '  +2.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '51.47'
$ws.Range("E32").Value = '  -1.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.29'
$ws.Range("E33").Value = '  -1.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0445'
$ws.Range("E34").Value = '  -8.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0834'
$ws.Range("E35").Value = '  -2.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.15'
$ws.Range("E36").Value = '  -8.27%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.91'
$ws.Range("E38").Value = '  +2.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.10'
$ws.Range("E39").Value = '  -5.35%  '
$ws.Range("E40").Value = '  -4.53%  '
$ws.Range("E41").Value = '  +1.77%  '
$ws.Range("E42").Value = '  -3.39%  '
$ws.Range("E43").Value = '  -2.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.50'
$ws.Range("E44").Value = '  -6.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.80'
$ws.Range("E45").Value = '  -6.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.085.86'
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.25'
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.27'
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.933'
$ws.Range("E49").Value = '  -4.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.54'
$ws.Range("E50").Value = '  -5.55%  '
$ws.Range("E51").Value = '  -6.55%  '
